$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999438518516,
    0.99709017844233783,
    0.9895527881575763,
    0.98984830389634781,
    0.98039632986798819,
    0.9575811705487014,
    0.95261143209800347,
    0.94642387762550151,
    0.93949388314879378,
    0.93348392126411617,
    0.93268248310159363,
    0.93141341970622071,
    0.92698104950410376,
    0.9257344390380231,
    0.92596498182823184,
    0.92684935662824697,
    0.92314155078905069,
    0.92203263839515071,
    0.99420970493439964,
    0.97086382694604667,
    0.96457532149167524,
    0.96331080314314654,
    0.98251832884524182,
    0.96949779795199276,
    0.96304085794251781,
    0.95553072176953058,
    0.95292193479129761,
    0.94189552148916322,
    0.93438831134622835,
    0.93166737409809453,
    0.93391530115165211,
    0.93568645443940712,
    0.93516644213797351
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
